$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.094.85'
$ws.Range("E2").Value = '  -1.73%  '

$ws.Range("D3").Value = '2.360.39'
$ws.Range("E3").Value = '  +2.10%  '

$ws.Range("E4").Value = '  +0.06%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '301.48'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.91%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '99.33'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.26%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.569'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.58%  '

$ws.Range("E8").Value = '  +0.04%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.512'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.90%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '34.43'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -3.67%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0799'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.03%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '7.14'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.49%  '

$ws.Range("E13").Value = '  -0.43%  '

$ws.Range("D14").Value = '2.723.88'
$ws.Range("E14").Value = '  +2.29%  '

$ws.Range("D15").Value = '2.364.47'
$ws.Range("E15").Value = '  +2.22%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.814'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.30%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '13.62'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -2.71%  '

$ws.Range("D18").Value = '46.022.03'
$ws.Range("E18").Value = '  -1.67%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '12.77'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -3.41%  '

$ws.Range("D20").Value = '0.0₃0966'
$ws.Range("E20").Value = '  +2.65%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.05'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.53%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '67.54'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +1.07%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '245.66'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.00%  '

$ws.Range("E24").Value = '  -2.48%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("E26").Value = '  -2.86%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '39.98'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -6.67%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.19'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.89%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.79'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.07%  '

$ws.Range("E30").Value = '  +21.62%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '20.99'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +4.00%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.78'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +6.28%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.52'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -4.11%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '146.42'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.50%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0776'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -2.88%  '

$ws.Range("E36").Value = '  -1.47%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.90'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +5.56%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.116'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.16%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '14.97'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -4.97%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '3.94'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -2.22%  '

$ws.Range("E41").Value = '  -2.02%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '3.22'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -5.57%  '

$ws.Range("D43").Value = '1.910.76'
$ws.Range("E43").Value = '  +3.69%  '

$ws.Range("E44").Value = '  +0.00%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '91.95'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.98%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.80'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -8.94%  '

$ws.Range("E47").Value = '  -6.26%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '8.33'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +4.74%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '98.06'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.92%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.592.49'
$ws.Range("E50").Value = '  +2.00%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '14.46'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +6.12%  '
